$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-08-13 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-08-14 Thursday", 2)

$d.Content.Find.Execute("799÷8=99, 7", $true, $false, $false, $false, $false, $true, 1, $false, "460÷9=51, 1", 2)
$d.Content.Find.Execute("521÷4=130, 1", $true, $false, $false, $false, $false, $true, 1, $false, "243÷2=121, 1", 2)
$d.Content.Find.Execute("767÷4=191, 3", $true, $false, $false, $false, $false, $true, 1, $false, "516÷3=172, 0", 2)
$d.Content.Find.Execute("904÷9=100, 4", $true, $false, $false, $false, $false, $true, 1, $false, "152÷5=30, 2", 2)
$d.Content.Find.Execute("726÷9=80, 6", $true, $false, $false, $false, $false, $true, 1, $false, "546÷2=273, 0", 2)

$d.Content.Find.Execute("661÷4=165, 1", $true, $false, $false, $false, $false, $true, 1, $false, "681÷4=170, 1", 2)
$d.Content.Find.Execute("995÷3=331, 2", $true, $false, $false, $false, $false, $true, 1, $false, "525÷5=105, 0", 2)
$d.Content.Find.Execute("123÷4=30, 3", $true, $false, $false, $false, $false, $true, 1, $false, "473÷5=94, 3", 2)
$d.Content.Find.Execute("853÷4=213, 1", $true, $false, $false, $false, $false, $true, 1, $false, "237÷2=118, 1", 2)
$d.Content.Find.Execute("448÷3=149, 1", $true, $false, $false, $false, $false, $true, 1, $false, "256÷3=85, 1", 2)

$d.Content.Find.Execute("303÷3=101, 0", $true, $false, $false, $false, $false, $true, 1, $false, "956÷9=106, 2", 2)
$d.Content.Find.Execute("208÷4=52, 0", $true, $false, $false, $false, $false, $true, 1, $false, "409÷5=81, 4", 2)
$d.Content.Find.Execute("426÷9=47, 3", $true, $false, $false, $false, $false, $true, 1, $false, "170÷3=56, 2", 2)
$d.Content.Find.Execute("132÷4=33, 0", $true, $false, $false, $false, $false, $true, 1, $false, "444÷3=148, 0", 2)
$d.Content.Find.Execute("401÷8=50, 1", $true, $false, $false, $false, $false, $true, 1, $false, "696÷7=99, 3", 2)

$d.Content.Find.Execute("490÷9=54, 4", $true, $false, $false, $false, $false, $true, 1, $false, "218÷6=36, 2", 2)
$d.Content.Find.Execute("273÷9=30, 3", $true, $false, $false, $false, $false, $true, 1, $false, "717÷7=102, 3", 2)
$d.Content.Find.Execute("859÷7=122, 5", $true, $false, $false, $false, $false, $true, 1, $false, "670÷7=95, 5", 2)
$d.Content.Find.Execute("467÷8=58, 3", $true, $false, $false, $false, $false, $true, 1, $false, "638÷8=79, 6", 2)
$d.Content.Find.Execute("854÷5=170, 4", $true, $false, $false, $false, $false, $true, 1, $false, "146÷4=36, 2", 2)

$d.Content.Find.Execute("731÷2=365, 1", $true, $false, $false, $false, $false, $true, 1, $false, "721÷6=120, 1", 2)
$d.Content.Find.Execute("941÷3=313, 2", $true, $false, $false, $false, $false, $true, 1, $false, "623÷9=69, 2", 2)
$d.Content.Find.Execute("961÷4=240, 1", $true, $false, $false, $false, $false, $true, 1, $false, "278÷2=139, 0", 2)
$d.Content.Find.Execute("380÷4=95, 0", $true, $false, $false, $false, $false, $true, 1, $false, "433÷8=54, 1", 2)
$d.Content.Find.Execute("698÷7=99, 5", $true, $false, $false, $false, $false, $true, 1, $false, "444÷4=111, 0", 2)
